$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Read the old values before we start overwriting (A:D, rows 1-9) ---
# Old layout: A=Distribución, B=K-S (p), C=A-D (stat), D=Tipo

# Capture old data (rows 2-9) so we can re-place it in new columns.
$oldDist = @{}
$oldKs   = @{}
$oldAd   = @{}
$oldTipo = @{}
for ($r = 2; $r -le 9; $r++) {
    $oldDist[$r] = $ws.Cells.Item($r, 1).Value()
    $oldKs[$r]   = $ws.Cells.Item($r, 2).Value()
    $oldAd[$r]   = $ws.Cells.Item($r, 3).Value()
    $oldTipo[$r] = $ws.Cells.Item($r, 4).Value()
}

# --- New headers (row 1): Tipo, Distribución, K-S (p), A-D (stat), Evaluación Visual, Decisión Final, Letra Notación ---
$ws.Cells.Item(1, 1).Value = "Tipo"
$ws.Cells.Item(1, 2).Value = "Distribución"
$ws.Cells.Item(1, 3).Value = "K-S (p)"
$ws.Cells.Item(1, 4).Value = "A-D (stat)"
$ws.Cells.Item(1, 5).Value = "Evaluación Visual"
$ws.Cells.Item(1, 6).Value = "Decisión Final"
$ws.Cells.Item(1, 7).Value = "Letra Notación"

# Mirror the header style (bold, centered, top-aligned, thin box border),
# already applied to A1:D1, onto the new header cells E1:G1.
$newHeaderRange = $ws.Range("E1:G1")
$newHeaderRange.Font.Bold = $true
$newHeaderRange.HorizontalAlignment = -4108  # xlCenter
$newHeaderRange.VerticalAlignment = -4160    # xlTop
$newHeaderRange.Borders.LineStyle = 1        # xlContinuous (thin box border)

# --- New data rows: shift old A->B, old B->C, old C->D, old D(Tipo)->A, then add E/F/G ---
for ($r = 2; $r -le 9; $r++) {
    $ws.Cells.Item($r, 1).Value = $oldTipo[$r]
    $ws.Cells.Item($r, 2).Value = $oldDist[$r]
    $ws.Cells.Item($r, 3).Value = $oldKs[$r]
    $ws.Cells.Item($r, 4).Value = $oldAd[$r]
}

# Relabel "Llegadas" -> "Tiempos entre llegadas" and "Servicio" -> "Tiempos de servicio"
$ws.Range("A2:A5").Value = "Tiempos entre llegadas"
$ws.Range("A6:A9").Value = "Tiempos de servicio"

# Slightly refined A-D stat values (recomputed with the new evaluation pass)
$ws.Cells.Item(3, 4).Value = 7.983629470097156
$ws.Cells.Item(7, 4).Value = 4.148560891450487

# Evaluación Visual (column E)
$ws.Cells.Item(2, 5).Value = "Aceptable"
$ws.Cells.Item(3, 5).Value = "Malo"
$ws.Cells.Item(4, 5).Value = "Aceptable"
$ws.Cells.Item(5, 5).Value = "Malo"
$ws.Cells.Item(6, 5).Value = "Malo"
$ws.Cells.Item(7, 5).Value = "Regular"
$ws.Cells.Item(8, 5).Value = "Excelente"
$ws.Cells.Item(9, 5).Value = "Bueno"

# Decisión Final (column F)
$ws.Cells.Item(2, 6).Value = "✅ Se acepta"
$ws.Cells.Item(3, 6).Value = "❌ Se descarta"
$ws.Cells.Item(4, 6).Value = "✅ Se acepta"
$ws.Cells.Item(5, 6).Value = "❌ Se descarta"
$ws.Cells.Item(6, 6).Value = "❌ Se descarta"
$ws.Cells.Item(7, 6).Value = "✅ Se acepta"
$ws.Cells.Item(8, 6).Value = "✅ Se acepta"
$ws.Cells.Item(9, 6).Value = "✅ Se acepta"

# Letra Notación (column G)
$ws.Cells.Item(2, 7).Value = "M"
$ws.Cells.Item(3, 7).Value = "N"
$ws.Cells.Item(4, 7).Value = "G"
$ws.Cells.Item(5, 7).Value = "G"
$ws.Cells.Item(6, 7).Value = "M"
$ws.Cells.Item(7, 7).Value = "N"
$ws.Cells.Item(8, 7).Value = "G"
$ws.Cells.Item(9, 7).Value = "G"
